$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 gets replaced with new login-test data
$ws.Range("A2").Value = "NinjaAlgo1"
$ws.Range("B2").Value = "@Algo3"
$ws.Range("C2").Value = "login"

# New row 3 holds the original credentials with an updated status
$ws.Range("A3").Value = "NinjaAlgo"
$ws.Range("B3").Value = "@Algo123"
$ws.Range("C3").Value = "home"
